# Update 22-Jan-2021, end of day update.
# Applies the petty-cash book edits described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Existing rows: add extra spend/income to already-present formulas ---
$ws.Range("D17").Formula = "=45000+195000"
$ws.Range("D18").Formula = "=14805000+2565000+336000+337000"
$ws.Range("C19").Formula = "=2580000+1410000+28052000"

# --- New rows 22-24 (still 20-Jan-2021 block) ---
$ws.Range("B22").Value = "SALES - cash/retail"
$ws.Range("C22").Formula = "=13661975+30077025-28052000"

$ws.Range("B23").Value = "SELISIH - lebih"
$ws.Range("C23").Value = 140000

$ws.Range("B24").Value = "SETOR KE BANK"
$ws.Range("D24").Value = 29000000

# --- New rows 25-33 (21-Jan-2021 block) ---
$ws.Range("A25").Value = 44217
$ws.Range("B25").Value = "Wages Expense"
$ws.Range("D25").Formula = "=45000+225000"

$ws.Range("B26").Value = "TRANSFER BCA"
$ws.Range("D26").Formula = "=1630000+2050000+280000+371500"

$ws.Range("B27").Value = "CHEQUE RECEIVED"
$ws.Range("D27").Formula = "=1619000"

$ws.Range("B28").Value = "A/P"
$ws.Range("D28").Formula = "=3925000"

$ws.Range("B29").Value = "A/R"
$ws.Range("C29").Formula = "=18910500"

$ws.Range("B30").Value = "SALES - cash/retail"
$ws.Range("C30").Formula = "=8526975+21744025-18910500"

$ws.Range("B31").Value = "prive"
$ws.Range("D31").Value = 1000000

$ws.Range("B32").Value = "SELISIH - lebih"
$ws.Range("C32").Value = 420000

$ws.Range("B33").Value = "SETOR KE BANK"
$ws.Range("D33").Value = 20000000

# --- New row 34 marks the start of 22-Jan-2021 ---
$ws.Range("A34").Value = 44218

# --- View state: selection moved further down the sheet as work progressed ---
$ws.Activate()
$ws.Range("C54").Select()

$wb.Application.CalculateFullRebuild()
